# feat: add 2022-Q1 data
#
# Before: 4 sheets -> 2021-Q2, 2021-Q3, 2021-Q4, 总计
# After:  5 sheets -> 2021-Q2, 2021-Q3, 2021-Q4, 2022-Q1, 总计
#
# A new "2022-Q1" fund-holdings sheet is inserted right before the
# "总计" (totals) sheet, and the totals sheet gets a new first data row
# summarizing the 2022-Q1 quarter (existing rows shift down by one).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: turn the existing "总计" sheet into the new "2022-Q1" sheet.
# Reusing this sheet object (instead of creating a brand new one) keeps
# its sheetId/sheetPr/pageMargins/etc. intact, and it ends up 4th in
# tab order right where "总计" used to sit.
# ---------------------------------------------------------------------
$q1 = $wb.Worksheets.Item("总计")
$q1.Name = "2022-Q1"
$q1.Cells.Clear() | Out-Null

# A same-formatted sheet to borrow header/index cell styling from.
$template = $wb.Worksheets.Item("2021-Q4")
$headerSrc = $template.Range("B1")
$idxSrc = $template.Range("A2")

$fundHeaders = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
$col = 2
foreach ($h in $fundHeaders) {
    $cell = $q1.Cells.Item(1, $col)
    $headerSrc.Copy() | Out-Null
    $cell.PasteSpecial(-4122) | Out-Null
    $cell.Value = $h
    $col = $col + 1
}

function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats() | Out-Null
}

$idxCell = $q1.Cells.Item(2, 1)
$idxSrc.Copy() | Out-Null
$idxCell.PasteSpecial(-4122) | Out-Null
$idxCell.Value = 0

Set-TextValue $q1.Cells.Item(2, 2) "006477"
Set-TextValue $q1.Cells.Item(2, 3) "中邮沪港深精选混合"
Set-TextValue $q1.Cells.Item(2, 4) "0.05"
Set-TextValue $q1.Cells.Item(2, 5) "83.91"
Set-TextValue $q1.Cells.Item(2, 6) "4.53"
Set-TextValue $q1.Cells.Item(2, 7) "0.0023"
$q1.Cells.Item(2, 8).Value = 7

# ---------------------------------------------------------------------
# Step 2: add a brand new "总计" sheet right after "2022-Q1" by copying
# a same-formatted sheet (so sheetPr/pageMargins/sheetFormatPr come out
# matching the rest of the workbook instead of Excel's blank defaults),
# then clearing its copied content and writing the updated totals.
# ---------------------------------------------------------------------
$template.Copy($null, $q1) | Out-Null
$total = $wb.Worksheets.Item(5)
$total.Name = "总计"
$total.Cells.Clear() | Out-Null

$totalsHeaders = @("日期", "持有数量(只)", "持有市值(亿元)")
$col = 2
foreach ($h in $totalsHeaders) {
    $cell = $total.Cells.Item(1, $col)
    $headerSrc.Copy() | Out-Null
    $cell.PasteSpecial(-4122) | Out-Null
    $cell.Value = $h
    $col = $col + 1
}

function Add-TotalsRow($rowNum, $idx, $date, $count, $value) {
    $idxCell = $total.Cells.Item($rowNum, 1)
    $idxSrc.Copy() | Out-Null
    $idxCell.PasteSpecial(-4122) | Out-Null
    $idxCell.Value = $idx

    $total.Cells.Item($rowNum, 2).Value = $date
    $total.Cells.Item($rowNum, 3).Value = $count
    $total.Cells.Item($rowNum, 4).Value = $value
}

Add-TotalsRow 2 0 "2022-Q1" 1 0
Add-TotalsRow 3 1 "2021-Q4" 4 0.1
Add-TotalsRow 4 2 "2021-Q3" 1 0.02
Add-TotalsRow 5 3 "2021-Q2" 5 1.27

# Restore the originally-active sheet/selection (unchanged by the diff).
$wb.Worksheets.Item("2021-Q2").Activate()
